# Fruta / hortaliza, semanal
# Insert a new daily price record as row 99 (Macroferia Regional de Talca - Mango),
# pushing all the following records down by one row. The last existing record
# (old row 152) ends up at row 153, matching the new dimension A1:T153.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 99..152 down to 100..153, leaving row 99 free for the new record.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new data point.
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44845
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100108
$ws.Range("H99").Value = "Tropicales y subtropicales"
$ws.Range("I99").Value = 100108002
$ws.Range("J99").Value = "Mango"
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "Primera"
$ws.Range("M99").Value = 300
$ws.Range("N99").Value = 7000
$ws.Range("O99").Value = 7000
$ws.Range("P99").Value = 7000
$ws.Range("Q99").Value = "$/bandeja 4 kilos"
$ws.Range("R99").Value = "Brasil"
$ws.Range("S99").Value = 1750
$ws.Range("T99").Value = 4
